$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 78
$ws.Range("E3").Value = 29
$ws.Range("E5").Value = 99
$ws.Range("E6").Value = 36
$ws.Range("E7").Value = 26
$ws.Range("E10").Value = 334
$ws.Range("E11").Value = 229
$ws.Range("E12").Value = 334
$ws.Range("F12").Value = 190
$ws.Range("H12").Value = 190
$ws.Range("E13").Value = 100
$ws.Range("F13").Value = 53
$ws.Range("H13").Value = 53
$ws.Range("E14").Value = 91
$ws.Range("E16").Value = 146
$ws.Range("F16").Value = 76
$ws.Range("H16").Value = 76
$ws.Range("E18").Value = 45
$ws.Range("E20").Value = 70
$ws.Range("E21").Value = 110
$ws.Range("F21").Value = 58
$ws.Range("H21").Value = 58
$ws.Range("E22").Value = 128
$ws.Range("E23").Value = 149
$ws.Range("F23").Value = 68
$ws.Range("H23").Value = 68
$ws.Range("E24").Value = 149
$ws.Range("E25").Value = 173
$ws.Range("E26").Value = 94
$ws.Range("F26").Value = 58
$ws.Range("H26").Value = 58
$ws.Range("E27").Value = 228
$ws.Range("E28").Value = 138
$ws.Range("F28").Value = 47
$ws.Range("H28").Value = 47
$ws.Range("E29").Value = 135
$ws.Range("F29").Value = 78
$ws.Range("H29").Value = 78
$ws.Range("E30").Value = 151
$ws.Range("F30").Value = 86
$ws.Range("H30").Value = 86
$ws.Range("E31").Value = 54
$ws.Range("E32").Value = 134
$ws.Range("F32").Value = 73
$ws.Range("H32").Value = 73
$ws.Range("E33").Value = 216
$ws.Range("F33").Value = 107
$ws.Range("H33").Value = 107
$ws.Range("E34").Value = 159
$ws.Range("E35").Value = 101
$ws.Range("E36").Value = 46
$ws.Range("E37").Value = 114
$ws.Range("E38").Value = 72
$ws.Range("F38").Value = 49
$ws.Range("H38").Value = 49
$ws.Range("E39").Value = 148
$ws.Range("F39").Value = 68
$ws.Range("H39").Value = 68
$ws.Range("E40").Value = 193
$ws.Range("E41").Value = 284
$ws.Range("E42").Value = 249
$ws.Range("F42").Value = 133
$ws.Range("H42").Value = 133
$ws.Range("E43").Value = 82
$ws.Range("E44").Value = 230
$ws.Range("E45").Value = 98
$ws.Range("E46").Value = 218
$ws.Range("F46").Value = 118
$ws.Range("H46").Value = 118
$ws.Range("E47").Value = 318
$ws.Range("F47").Value = 156
$ws.Range("H47").Value = 156
$ws.Range("E48").Value = 144
$ws.Range("E49").Value = 206
$ws.Range("F49").Value = 88
$ws.Range("H49").Value = 88
$ws.Range("E50").Value = 177
$ws.Range("E51").Value = 163
$ws.Range("F51").Value = 65
$ws.Range("H51").Value = 65
$ws.Range("E52").Value = 21
$ws.Range("F52").Value = 8
$ws.Range("H52").Value = 8
